$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: human readable labels
$ws.Range("A1").Value = "Horas trabajadas"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Situación profesional código"
$ws.Range("D1").Value = "Aragón"
$ws.Range("E1").Value = "Situación profesional"

# Row 2: concept URIs
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:situacion-profesional"

# Row 3: concept type (medida/dim)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"

# Row 4: datatype
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "xsd:string"
